$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData1")
$ws.Activate()

# Replace the text "Sup" in B4 with the numeric value 90
$ws.Range("B4").Value = 90

# Update the selection to reflect the newly edited cell
$ws.Range("B4").Select()
